$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-05 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-06 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("83-59=24", $true, $false, $false, $false, $false, $true, 1, $false, "26-5=21", 2) | Out-Null
$d.Content.Find.Execute("58-44=14", $true, $false, $false, $false, $false, $true, 1, $false, "47-23=24", 2) | Out-Null
$d.Content.Find.Execute("47+46=93", $true, $false, $false, $false, $false, $true, 1, $false, "92-20=72", 2) | Out-Null
$d.Content.Find.Execute("3+90=93", $true, $false, $false, $false, $false, $true, 1, $false, "3+40=43", 2) | Out-Null
$d.Content.Find.Execute("19+0=19", $true, $false, $false, $false, $false, $true, 1, $false, "94-90=4", 2) | Out-Null
$d.Content.Find.Execute("67-50=17", $true, $false, $false, $false, $false, $true, 1, $false, "96-24=72", 2) | Out-Null
$d.Content.Find.Execute("3+23=26", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=7", 2) | Out-Null
$d.Content.Find.Execute("64+1=65", $true, $false, $false, $false, $false, $true, 1, $false, "79+4=83", 2) | Out-Null
$d.Content.Find.Execute("62-11=51", $true, $false, $false, $false, $false, $true, 1, $false, "35+58=93", 2) | Out-Null
$d.Content.Find.Execute("31+12=43", $true, $false, $false, $false, $false, $true, 1, $false, "6+27=33", 2) | Out-Null
$d.Content.Find.Execute("54-1=53", $true, $false, $false, $false, $false, $true, 1, $false, "13+43=56", 2) | Out-Null
$d.Content.Find.Execute("46-34=12", $true, $false, $false, $false, $false, $true, 1, $false, "23+53=76", 2) | Out-Null
$d.Content.Find.Execute("36+8=44", $true, $false, $false, $false, $false, $true, 1, $false, "81-53=28", 2) | Out-Null
$d.Content.Find.Execute("1+43=44", $true, $false, $false, $false, $false, $true, 1, $false, "32+30=62", 2) | Out-Null
$d.Content.Find.Execute("55+7=62", $true, $false, $false, $false, $false, $true, 1, $false, "4+34=38", 2) | Out-Null
$d.Content.Find.Execute("72-71=1", $true, $false, $false, $false, $false, $true, 1, $false, "76+2=78", 2) | Out-Null
$d.Content.Find.Execute("69-48=21", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("26-9=17", $true, $false, $false, $false, $false, $true, 1, $false, "48+46=94", 2) | Out-Null
$d.Content.Find.Execute("32+3=35", $true, $false, $false, $false, $false, $true, 1, $false, "12+43=55", 2) | Out-Null
$d.Content.Find.Execute("9+67=76", $true, $false, $false, $false, $false, $true, 1, $false, "74-17=57", 2) | Out-Null
$d.Content.Find.Execute("23+24=47", $true, $false, $false, $false, $false, $true, 1, $false, "79-74=5", 2) | Out-Null
$d.Content.Find.Execute("98-43=55", $true, $false, $false, $false, $false, $true, 1, $false, "38+39=77", 2) | Out-Null
$d.Content.Find.Execute("0+58=58", $true, $false, $false, $false, $false, $true, 1, $false, "90-31=59", 2) | Out-Null
$d.Content.Find.Execute("9+60=69", $true, $false, $false, $false, $false, $true, 1, $false, "74-45=29", 2) | Out-Null
$d.Content.Find.Execute("73-38=35", $true, $false, $false, $false, $false, $true, 1, $false, "72-11=61", 2) | Out-Null
$d.Content.Find.Execute("99-26=73", $true, $false, $false, $false, $false, $true, 1, $false, "17+73=90", 2) | Out-Null
$d.Content.Find.Execute("67+31=98", $true, $false, $false, $false, $false, $true, 1, $false, "96-28=68", 2) | Out-Null
$d.Content.Find.Execute("58+4=62", $true, $false, $false, $false, $false, $true, 1, $false, "49-7=42", 2) | Out-Null
$d.Content.Find.Execute("51+35=86", $true, $false, $false, $false, $false, $true, 1, $false, "1+50=51", 2) | Out-Null
$d.Content.Find.Execute("59+14=73", $true, $false, $false, $false, $false, $true, 1, $false, "97-92=5", 2) | Out-Null
$d.Content.Find.Execute("76-46=30", $true, $false, $false, $false, $false, $true, 1, $false, "15-5=10", 2) | Out-Null
$d.Content.Find.Execute("40+30=70", $true, $false, $false, $false, $false, $true, 1, $false, "12+22=34", 2) | Out-Null
$d.Content.Find.Execute("5+76=81", $true, $false, $false, $false, $false, $true, 1, $false, "45+39=84", 2) | Out-Null
$d.Content.Find.Execute("56-55=1", $true, $false, $false, $false, $false, $true, 1, $false, "64-55=9", 2) | Out-Null
$d.Content.Find.Execute("49+47=96", $true, $false, $false, $false, $false, $true, 1, $false, "20+63=83", 2) | Out-Null
$d.Content.Find.Execute("52+20=72", $true, $false, $false, $false, $false, $true, 1, $false, "15+11=26", 2) | Out-Null
$d.Content.Find.Execute("29+27=56", $true, $false, $false, $false, $false, $true, 1, $false, "88-12=76", 2) | Out-Null
$d.Content.Find.Execute("24+39=63", $true, $false, $false, $false, $false, $true, 1, $false, "77-57=20", 2) | Out-Null
$d.Content.Find.Execute("46-20=26", $true, $false, $false, $false, $false, $true, 1, $false, "39-9=30", 2) | Out-Null
$d.Content.Find.Execute("69-1=68", $true, $false, $false, $false, $false, $true, 1, $false, "50+31=81", 2) | Out-Null
$d.Content.Find.Execute("32+10=42", $true, $false, $false, $false, $false, $true, 1, $false, "6+86=92", 2) | Out-Null
$d.Content.Find.Execute("96-43=53", $true, $false, $false, $false, $false, $true, 1, $false, "79-6=73", 2) | Out-Null
$d.Content.Find.Execute("90-14=76", $true, $false, $false, $false, $false, $true, 1, $false, "54-33=21", 2) | Out-Null
$d.Content.Find.Execute("78-32=46", $true, $false, $false, $false, $false, $true, 1, $false, "11+47=58", 2) | Out-Null
$d.Content.Find.Execute("9+42=51", $true, $false, $false, $false, $false, $true, 1, $false, "10-4=6", 2) | Out-Null
$d.Content.Find.Execute("43-43=0", $true, $false, $false, $false, $false, $true, 1, $false, "22+5=27", 2) | Out-Null
$d.Content.Find.Execute("98-54=44", $true, $false, $false, $false, $false, $true, 1, $false, "17-9=8", 2) | Out-Null
$d.Content.Find.Execute("77+5=82", $true, $false, $false, $false, $false, $true, 1, $false, "53+5=58", 2) | Out-Null
$d.Content.Find.Execute("49-35=14", $true, $false, $false, $false, $false, $true, 1, $false, "87-74=13", 2) | Out-Null
$d.Content.Find.Execute("94-72=22", $true, $false, $false, $false, $false, $true, 1, $false, "51-34=17", 2) | Out-Null
$d.Content.Find.Execute("44-34=10", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=18", 2) | Out-Null
$d.Content.Find.Execute("0+89=89", $true, $false, $false, $false, $false, $true, 1, $false, "73-14=59", 2) | Out-Null
$d.Content.Find.Execute("45+48=93", $true, $false, $false, $false, $false, $true, 1, $false, "32-15=17", 2) | Out-Null
$d.Content.Find.Execute("42-27=15", $true, $false, $false, $false, $false, $true, 1, $false, "3+46=49", 2) | Out-Null
$d.Content.Find.Execute("39+21=60", $true, $false, $false, $false, $false, $true, 1, $false, "32-27=5", 2) | Out-Null
$d.Content.Find.Execute("0+92=92", $true, $false, $false, $false, $false, $true, 1, $false, "74+6=80", 2) | Out-Null
$d.Content.Find.Execute("42+2=44", $true, $false, $false, $false, $false, $true, 1, $false, "16+61=77", 2) | Out-Null
$d.Content.Find.Execute("97-34=63", $true, $false, $false, $false, $false, $true, 1, $false, "72+12=84", 2) | Out-Null
$d.Content.Find.Execute("75-62=13", $true, $false, $false, $false, $false, $true, 1, $false, "46-38=8", 2) | Out-Null
$d.Content.Find.Execute("29+12=41", $true, $false, $false, $false, $false, $true, 1, $false, "91-31=60", 2) | Out-Null
$d.Content.Find.Execute("7+61=68", $true, $false, $false, $false, $false, $true, 1, $false, "57-9=48", 2) | Out-Null
$d.Content.Find.Execute("61-9=52", $true, $false, $false, $false, $false, $true, 1, $false, "73+3=76", 2) | Out-Null
$d.Content.Find.Execute("17-5=12", $true, $false, $false, $false, $false, $true, 1, $false, "8+22=30", 2) | Out-Null
$d.Content.Find.Execute("48-1=47", $true, $false, $false, $false, $false, $true, 1, $false, "62+9=71", 2) | Out-Null
$d.Content.Find.Execute("66+29=95", $true, $false, $false, $false, $false, $true, 1, $false, "90+4=94", 2) | Out-Null
$d.Content.Find.Execute("26+22=48", $true, $false, $false, $false, $false, $true, 1, $false, "37+58=95", 2) | Out-Null
$d.Content.Find.Execute("29+56=85", $true, $false, $false, $false, $false, $true, 1, $false, "67+0=67", 2) | Out-Null
$d.Content.Find.Execute("92-14=78", $true, $false, $false, $false, $false, $true, 1, $false, "65-35=30", 2) | Out-Null
$d.Content.Find.Execute("26+57=83", $true, $false, $false, $false, $false, $true, 1, $false, "21-1=20", 2) | Out-Null
$d.Content.Find.Execute("6+41=47", $true, $false, $false, $false, $false, $true, 1, $false, "4+36=40", 2) | Out-Null
$d.Content.Find.Execute("56+32=88", $true, $false, $false, $false, $false, $true, 1, $false, "82+13=95", 2) | Out-Null
$d.Content.Find.Execute("52-6=46", $true, $false, $false, $false, $false, $true, 1, $false, "39-18=21", 2) | Out-Null
$d.Content.Find.Execute("83-81=2", $true, $false, $false, $false, $false, $true, 1, $false, "77-30=47", 2) | Out-Null
$d.Content.Find.Execute("91-78=13", $true, $false, $false, $false, $false, $true, 1, $false, "75+12=87", 2) | Out-Null
$d.Content.Find.Execute("17+26=43", $true, $false, $false, $false, $false, $true, 1, $false, "36+13=49", 2) | Out-Null
$d.Content.Find.Execute("62+34=96", $true, $false, $false, $false, $false, $true, 1, $false, "86-30=56", 2) | Out-Null
$d.Content.Find.Execute("66+10=76", $true, $false, $false, $false, $false, $true, 1, $false, "80+7=87", 2) | Out-Null
$d.Content.Find.Execute("90-65=25", $true, $false, $false, $false, $false, $true, 1, $false, "7+16=23", 2) | Out-Null
$d.Content.Find.Execute("20+43=63", $true, $false, $false, $false, $false, $true, 1, $false, "84-7=77", 2) | Out-Null
$d.Content.Find.Execute("21+11=32", $true, $false, $false, $false, $false, $true, 1, $false, "54+44=98", 2) | Out-Null
$d.Content.Find.Execute("90+6=96", $true, $false, $false, $false, $false, $true, 1, $false, "87-78=9", 2) | Out-Null
$d.Content.Find.Execute("30+53=83", $true, $false, $false, $false, $false, $true, 1, $false, "17+22=39", 2) | Out-Null
$d.Content.Find.Execute("81-75=6", $true, $false, $false, $false, $false, $true, 1, $false, "93-88=5", 2) | Out-Null
$d.Content.Find.Execute("46+16=62", $true, $false, $false, $false, $false, $true, 1, $false, "55-21=34", 2) | Out-Null
$d.Content.Find.Execute("26+5=31", $true, $false, $false, $false, $false, $true, 1, $false, "41-32=9", 2) | Out-Null
$d.Content.Find.Execute("93-72=21", $true, $false, $false, $false, $false, $true, 1, $false, "46-21=25", 2) | Out-Null
$d.Content.Find.Execute("16+39=55", $true, $false, $false, $false, $false, $true, 1, $false, "82-78=4", 2) | Out-Null
$d.Content.Find.Execute("11+34=45", $true, $false, $false, $false, $false, $true, 1, $false, "14+3=17", 2) | Out-Null
$d.Content.Find.Execute("50-10=40", $true, $false, $false, $false, $false, $true, 1, $false, "7+87=94", 2) | Out-Null
$d.Content.Find.Execute("56-13=43", $true, $false, $false, $false, $false, $true, 1, $false, "37-6=31", 2) | Out-Null
$d.Content.Find.Execute("83-32=51", $true, $false, $false, $false, $false, $true, 1, $false, "35+14=49", 2) | Out-Null
$d.Content.Find.Execute("22+21=43", $true, $false, $false, $false, $false, $true, 1, $false, "48-21=27", 2) | Out-Null
$d.Content.Find.Execute("47+20=67", $true, $false, $false, $false, $false, $true, 1, $false, "31+32=63", 2) | Out-Null
$d.Content.Find.Execute("90-49=41", $true, $false, $false, $false, $false, $true, 1, $false, "17+52=69", 2) | Out-Null
$d.Content.Find.Execute("86+12=98", $true, $false, $false, $false, $false, $true, 1, $false, "2+54=56", 2) | Out-Null
$d.Content.Find.Execute("85-58=27", $true, $false, $false, $false, $false, $true, 1, $false, "62-55=7", 2) | Out-Null
$d.Content.Find.Execute("65-24=41", $true, $false, $false, $false, $false, $true, 1, $false, "79-37=42", 2) | Out-Null
$d.Content.Find.Execute("21+62=83", $true, $false, $false, $false, $false, $true, 1, $false, "46+0=46", 2) | Out-Null
$d.Content.Find.Execute("75-47=28", $true, $false, $false, $false, $false, $true, 1, $false, "54+27=81", 2) | Out-Null
$d.Content.Find.Execute("0+96=96", $true, $false, $false, $false, $false, $true, 1, $false, "32-10=22", 2) | Out-Null

Write-Output "Replacements complete: 101 items"
